# Add a new parameter row "general.maxNumberCompThreads" to the
# structuralParameterProperties sheet, right before the existing
# "general.maxMemoryGB" row (old row 21), pushing it and every row
# below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 21; everything currently at/after row 21
# (general.maxMemoryGB, ...) shifts down to row 22, etc.
$ws.Rows.Item(21).Insert()

# Columns: A=name, B=inputFile, C=outputFile, D=parameter, E=type,
# F=attributes, G=usage, H=description
$ws.Range("A21").Value = "general.maxNumberCompThreads"
$ws.Range("D21").Value = "structural_pipeline"
$ws.Range("E21").Value = "numeric"
$ws.Range("F21").Value = "scalar nonempty nonnegative"
$ws.Range("G21").Value = "standard"
$ws.Range("H21").Value = "Maximum number of computational threads used in pipeline. Value 0 lets MATLAB determine the most desirable number of computational threads (equal to the number of physical cores on the machine)."

# Match the styling used by neighboring "standard" rows (numFmt text
# style on the attributes/usage columns).
$ws.Range("F21").NumberFormat = $ws.Range("F22").NumberFormat
$ws.Range("G21").NumberFormat = $ws.Range("G22").NumberFormat

$ws.Range("H22").Select()
